# edit.ps1 — reproduce the commit's changes against before.pptx
#
# 1) The table on slide 6 gets its table style (tableStyleId) switched from
#    the custom "Table_0" style ({503BFA9C-E17C-4AD3-8A2D-4B5E7FF54F62}) to
#    the built-in style {A1316A0F-16C7-4F38-971A-D07E6B0687C7}.
#
# 2) ppt/theme/theme1.xml and ppt/theme/theme2.xml swap their content
#    ("Office Theme" <-> "Integral"). theme2.xml is the theme actually used
#    by the (only) slide master / all slides, so we recolor it to the
#    "Office Theme" palette via the ThemeColorScheme object, which is the
#    piece of that swap reachable from the PowerPoint object model.

$p = $ppt.ActivePresentation

# --- 1) Table style on slide 6 -------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{A1316A0F-16C7-4F38-971A-D07E6B0687C7}")
    }
}

# --- 2) Theme colours: Integral -> Office Theme ---------------------------
function Set-ThemeRGB {
    param($scheme, [int]$index, [int]$r, [int]$g, [int]$b)
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$colorScheme = $p.Slides.Item(1).ThemeColorScheme

Set-ThemeRGB $colorScheme 1  0   0   0      # dk1
Set-ThemeRGB $colorScheme 2  255 255 255    # lt1
Set-ThemeRGB $colorScheme 3  68  84  106    # dk2
Set-ThemeRGB $colorScheme 4  231 230 230    # lt2
Set-ThemeRGB $colorScheme 5  91  155 213    # accent1
Set-ThemeRGB $colorScheme 6  237 125 49     # accent2
Set-ThemeRGB $colorScheme 7  165 165 165    # accent3
Set-ThemeRGB $colorScheme 8  255 192 0      # accent4
Set-ThemeRGB $colorScheme 9  68  114 196    # accent5
Set-ThemeRGB $colorScheme 10 112 173 71     # accent6
Set-ThemeRGB $colorScheme 11 5   99  193    # hlink
Set-ThemeRGB $colorScheme 12 149 79  114    # folHlink
